# Small changes to existing phrases on the "list" sheet:
#   - "Ducks in a row"          -> "ducks in a row"
#   - "Chewing sounds"          -> "chewing sounds"
#   - "Animal or child sounds"  -> split into "animal sounds" (A53) and a
#     new trailing row "child sounds" (A57)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$ws.Range("A51").Value = "ducks in a row"
$ws.Range("A52").Value = "chewing sounds"
$ws.Range("A53").Value = "animal sounds"
$ws.Range("A57").Value = "child sounds"

$ws.Activate() | Out-Null
$ws.Range("A56").Select() | Out-Null
